# Generate Report for Handoff
#
# The "b.md" file has moved from "Handed back: in sync with en-US" to a
# new handoff ("Ready for handoff") with a freshly generated handoff
# file (b.63290e5768f688058c7b37413b0a5c26c308f864.<locale>.xlf) and an
# updated handoff datetime. Update the Overview sheet plus the per-locale
# (zh-cn / de-de) detail sheets, including the "Latest Handoff File"
# hyperlink, which must now point its *display text* at the new xlf file
# name while still resolving to the same handoff-commit URL used before.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is "b.md" -> zh-cn / de-de status + handoff date
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value2 = "Ready for handoff"
$ov.Range("C3").Value2 = "Ready for handoff"
$ov.Range("D3").Value2 = "2016-26-12 16:26:53"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 ("b.md") Status / Latest Handoff File / Datetime
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value2 = "Ready for handoff"
$zh.Range("D3").Value2 = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("E3").Value2 = "2016-03-12 16:26:50"

# Rebuild the hyperlinks for this sheet so the "Latest Handoff File"
# hyperlink display text on D3 reflects the new file name while every
# other hyperlink (and every target URL, including D3's own) is kept
# exactly as it was.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ad294fb2753af7d136a220b6d2c54a0c79afccab/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/ad294fb2753af7d136a220b6d2c54a0c79afccab/e2e/a.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da314df305d29918eafea66652b186d53a2627d7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e0d8b010a069d3b0cdc52f3ad7362efd0d0ba8e9/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a59ddbce76ffdbbe93df6fd4003ff1ab066c4e12/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ad294fb2753af7d136a220b6d2c54a0c79afccab/e2e/b.md", "", "", "b.md")
$zh.Hyperlinks.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/ad294fb2753af7d136a220b6d2c54a0c79afccab/e2e/b.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da314df305d29918eafea66652b186d53a2627d7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e0d8b010a069d3b0cdc52f3ad7362efd0d0ba8e9/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a59ddbce76ffdbbe93df6fd4003ff1ab066c4e12/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet: row 3 ("b.md") Status / Latest Handoff File / Datetime
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value2 = "Ready for handoff"
$de.Range("D3").Value2 = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("E3").Value2 = "2016-03-12 16:26:53"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ad294fb2753af7d136a220b6d2c54a0c79afccab/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/ad294fb2753af7d136a220b6d2c54a0c79afccab/e2e/a.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/24890f2cdcd59ea060be313e37817dc35f4afe21/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8b6a13fd54870a66d5e78ad9b715eef8a82c6bc7/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1995f0963f82a219a38f755683b6154de423a257/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ad294fb2753af7d136a220b6d2c54a0c79afccab/e2e/b.md", "", "", "b.md")
$de.Hyperlinks.Add($de.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/ad294fb2753af7d136a220b6d2c54a0c79afccab/e2e/b.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/24890f2cdcd59ea060be313e37817dc35f4afe21/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8b6a13fd54870a66d5e78ad9b715eef8a82c6bc7/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1995f0963f82a219a38f755683b6154de423a257/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")

Write-Output "Generated handoff report update for b.md (zh-cn, de-de)."
